$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (column C) date value for every existing data row
#    (rows 2-394) from 45190 to 45192.
$ws.Range("C2:C394").Value2 = 45192

# 2. Rows 393 and 394 swap their "Beteckning" (A) and "Area (ha)" (G) values.
#    Row 393 becomes what row 394 used to be, and vice versa. Row 394 also
#    gains an explicit row height (matching the other data rows).
$ws.Range("A393").Value2 = "A 44596-2023"
$ws.Range("G393").Value2 = 4.1

$ws.Range("A394").Value2 = "A 44495-2023"
$ws.Range("G394").Value2 = 1.1
$ws.Range("A394").EntireRow.RowHeight = 15

# 3. Append six new data rows (395-400) with the same shape as the existing
#    rows (Beteckning, Datum, Förändrad, Län, Kommun, Area, and the zeroed
#    species-count columns H-Q, plus a blank styled "Artnamn" cell in R).
$newRows = @(
    @{ Row = 395; A = "A 45102-2023"; G = 0.4 },
    @{ Row = 396; A = "A 44910-2023"; G = 5.3 },
    @{ Row = 397; A = "A 45103-2023"; G = 3.4 },
    @{ Row = 398; A = "A 45101-2023"; G = 0.8 },
    @{ Row = 399; A = "A 44913-2023"; G = 11.7 },
    @{ Row = 400; A = "A 44929-2023"; G = 13.9 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value2 = $r.A                     # A - Beteckning
    $ws.Cells.Item($row, 2).Value2 = 45190                    # B - Datum
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($row, 3).Value2 = 45192                    # C - Förändrad
    $ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($row, 4).Value2 = "VÄSTERBOTTENS LÄN"       # D - Län
    $ws.Cells.Item($row, 5).Value2 = "MALÅ"                   # E - Kommun
    $ws.Cells.Item($row, 7).Value2 = $r.G                     # G - Area (ha)

    for ($col = 8; $col -le 17; $col++) {                     # H..Q = 0
        $ws.Cells.Item($row, $col).Value2 = 0
    }

    $ws.Cells.Item($row, 18).WrapText = $true                 # R - Artnamn (blank, wrapped)

    if ($row -le 399) {
        $ws.Cells.Item($row, 1).EntireRow.RowHeight = 15
    }
}
